$d = $word.ActiveDocument

# Anchor on the last paragraph in the document ("$sudo apt install ... ")
# and insert the new block of paragraphs right before its paragraph mark,
# so the new paragraphs do NOT inherit the ListParagraph/numbering
# formatting of that anchor paragraph.
$anchor = $d.Paragraphs.Last.Range
$insertAt = $anchor.End - 1
$target = $d.Range($insertAt, $insertAt)

$newParagraphsXml = @'
<w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>$sudo vi /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/postgresql/11/main/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postgresql.conf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">            Ensure   </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>listening_address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = “*”</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t xml:space="preserve">                           Port = 5433</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:t>$sudo vi /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/postgresql/11/main/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pg_hba.conf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Launch pgAdmin4 </w:t></w:r></w:p>
'@

$target.InsertXML($newParagraphsXml)
